$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Periodo Mora" for MIRELIS VALIENTE MARTINEZ (row 16) and
# EDINSON VELASQUEZ PALACIN (row 17) from 2506 to 2508. Both rows share
# the same period text, so update both cells to keep them in sync.
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"

# Update "Salario Basico" for MIRELIS VALIENTE MARTINEZ (row 16) with the
# new account statement value.
$ws.Range("G16").Value = 1423500
